# Update control flow ppt slides: convert the C# console-app snippets
# ("using static System.Console;", "WriteLine", "Write", "ReadLine", ...)
# into SplashKit-style C snippets ("#include \"splashkit.h\"", "write_line",
# "write", "read_line", ...).
#
# Helper: replace the text of a TextRange/Characters sub-range with brand
# new content without the engine's prefix/suffix auto-diff silently
# re-splitting the run in two. We first stomp the range with a same-length
# placeholder that shares no characters with either the old or the new
# text, then assign the real text - that second assignment has no common
# prefix/suffix with the placeholder, so it lands as a single clean run.
function Set-RangeText($range, $newText) {
    $placeholder = "@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@@"
    $len = $range.Text.Length
    $ph = $placeholder
    while ($ph.Length -lt $len) {
        $ph = $ph + $placeholder
    }
    if ($len -gt 0) {
        $ph = $ph.Substring(0, $len)
        $range.Text = $ph
    }
    $range.Text = $newText
}

# Locate $oldSub inside the (current) text of paragraph $para and replace
# just that substring's run with $newText, leaving the rest of the
# paragraph's runs untouched.
function Replace-Substring-In-Paragraph($tr, $para, $oldSub, $newText) {
    $ptext = $para.Text
    $idx = $ptext.IndexOf($oldSub)
    if ($idx -lt 0) {
        throw ("Substring not found: [" + $oldSub + "] in [" + $ptext + "]")
    }
    $absStart = $para.Start + $idx
    $sub = $tr.Characters($absStart, $oldSub.Length)
    Set-RangeText $sub $newText
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: "TextBox 15" inside "Group 11" - the using-directives / intro
# WriteLine calls.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$g1 = $s1.Shapes.Item(1)
$tb1 = $g1.GroupItems.Item(1)
$tr1 = $tb1.TextFrame.TextRange

# Paragraph 1: "using static System.Console;" -> "// include SplashKit library"
$para = $tr1.Paragraphs(1)
$newMid = '// include SplashKit library'
Replace-Substring-In-Paragraph $tr1 $para 'System.Console' $newMid
Replace-Substring-In-Paragraph $tr1 $para 'using static ' ''
Replace-Substring-In-Paragraph $tr1 $para ';' ''

# Paragraph 2: "using static System.Convert;" -> "#include "splashkit.h""
$para = $tr1.Paragraphs(2)
$newMid = '#include "splashkit.h"'
Replace-Substring-In-Paragraph $tr1 $para 'System.Convert' $newMid
Replace-Substring-In-Paragraph $tr1 $para 'using static ' ''
Replace-Substring-In-Paragraph $tr1 $para ';' ''

# Paragraph 4: WriteLine("Before you stands a 12 foot tall Knight...");
$para = $tr1.Paragraphs(4)
Set-RangeText $para 'write_line("Before you stands a 12 foot tall Knight...");'

# Paragraph 5: WriteLine();
$para = $tr1.Paragraphs(5)
Set-RangeText $para 'write_line();'

# Paragraph 6: WriteLine("\"We are the Knights who say 'Ni'.\"");
$para = $tr1.Paragraphs(6)
Set-RangeText $para 'write_line("\"We are the Knights who say ''Ni''.\"");'

# Paragraph 7: WriteLine("\"I will say Ni to you again if you do not appease us!\"");
$para = $tr1.Paragraphs(7)
Set-RangeText $para 'write_line("\"I will say Ni to you again if you do not appease us!\"");'

# ---------------------------------------------------------------------
# Slides 2-9: "TextBox 15" inside "Group 11" - the while(true) loop with
# the Knights-who-say-Ni sketch.
# ---------------------------------------------------------------------
for ($si = 2; $si -le 9; $si++) {
    $s = $p.Slides.Item($si)
    $g = $s.Shapes.Item(1)
    $tb = $g.GroupItems.Item(1)
    $tr = $tb.TextFrame.TextRange

    # Paragraph 3:       WriteLine("\"Ni!\"");
    $para = $tr.Paragraphs(3)
    Set-RangeText $para '      write_line("\"Ni!\"");'

    # Paragraph 4:       Write("Submit? ");
    $para = $tr.Paragraphs(4)
    Set-RangeText $para '      write("Submit? ");'

    # Paragraph 5:       if (ReadLine() == "y")  -- only the ReadLine run changes
    $para = $tr.Paragraphs(5)
    Replace-Substring-In-Paragraph $tr $para 'ReadLine' 'read_line'

    # Paragraph 11: WriteLine("\"Bring us a Shrubbery!\"");
    $para = $tr.Paragraphs(11)
    Set-RangeText $para 'write_line("\"Bring us a Shrubbery!\"");'
}

Write-Output "Done updating control-flow slides."
